$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-10 Thursday", "2025-04-11 Friday"),
    @("959÷8=119, 7", "452÷9=50, 2"),
    @("133÷2=66, 1", "447÷6=74, 3"),
    @("371÷2=185, 1", "963÷5=192, 3"),
    @("146÷2=73, 0", "578÷6=96, 2"),
    @("680÷2=340, 0", "745÷2=372, 1"),
    @("215÷4=53, 3", "639÷7=91, 2"),
    @("600÷6=100, 0", "728÷4=182, 0"),
    @("394÷8=49, 2", "806÷3=268, 2"),
    @("267÷7=38, 1", "957÷8=119, 5"),
    @("951÷9=105, 6", "747÷9=83, 0"),
    @("810÷7=115, 5", "231÷4=57, 3"),
    @("484÷4=121, 0", "353÷2=176, 1"),
    @("780÷6=130, 0", "998÷8=124, 6"),
    @("895÷8=111, 7", "995÷6=165, 5"),
    @("232÷3=77, 1", "583÷6=97, 1"),
    @("159÷9=17, 6", "455÷4=113, 3"),
    @("114÷2=57, 0", "950÷3=316, 2"),
    @("630÷5=126, 0", "745÷6=124, 1"),
    @("682÷7=97, 3", "253÷9=28, 1"),
    @("861÷5=172, 1", "755÷7=107, 6"),
    @("454÷3=151, 1", "908÷2=454, 0"),
    @("539÷3=179, 2", "690÷4=172, 2"),
    @("178÷9=19, 7", "103÷9=11, 4"),
    @("246÷8=30, 6", "918÷4=229, 2"),
    @("265÷3=88, 1", "860÷5=172, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
